# Added code for Create New Message functionality
# Adds a new "Messages" worksheet (subject/message) after the "Alerts" sheet,
# matching the freeCrmTestData.xlsx fixture used by the CRM QA test suite.

$wb = $excel.ActiveWorkbook

# --- Leftover UI state on the "Alerts" sheet (active prior to this edit) ---
$alerts = $wb.Worksheets.Item("Alerts")
$alerts.Activate()
$alerts.Columns.Item(4).AutoFit()
$alerts.Range("E14").Select()

# --- Create the new "Messages" sheet as the last tab (after "Alerts") -------
$messages = $wb.Worksheets.Add($null, $alerts)
$messages.Name = "Messages"

# --- Header row --------------------------------------------------------
$messages.Range("A1").Value = "subject"
$messages.Range("B1").Value = "message"
$messages.Range("A1:B1").Interior.Color = 65535

# --- Sample data rows ----------------------------------------------------
$messages.Range("A2").Value = "new message 1"
$messages.Range("B2").Value = "new message 1 text"

$messages.Range("A3").Value = "new message 2"
$messages.Range("B3").Value = "new message 2 text"

$messages.Range("A1").Select()
